$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6:C21").Select()
$ws.Range("A6:C21").ClearContents()
